$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "保險" (insurance) -> Worksheets.Item(6)
# Fix the header row (it previously duplicated row-2 data) and append the
# trailing metadata columns (property_category..index) to every row.
# ---------------------------------------------------------------------------
$wsIns = $wb.Worksheets.Item(6)

$insHeader = @{
    2  = "company"
    3  = "name"
    4  = "owner"
    5  = "property_category"
    6  = "category"
    7  = "date"
    8  = "legislator_name"
    9  = "legislator_id"
    10 = "source_file"
    11 = "index"
}
foreach ($col in $insHeader.Keys) {
    $wsIns.Cells.Item(1, $col).Value = $insHeader[$col]
}

$insRow2 = @{
    2  = "富邦人壽"
    3  = "投資型壽險"
    4  = "姚文智"
    5  = "insurance"
    6  = "normal"
    7  = "2012-04-30"
    8  = "姚文智"
    9  = 1745
    10 = "tmp60da1"
    11 = 76
}
foreach ($col in $insRow2.Keys) {
    $wsIns.Cells.Item(2, $col).Value = $insRow2[$col]
}

$insRow3 = @{
    2  = "南山人壽"
    3  = "投資型壽險"
    4  = "潘瓊琪"
    5  = "insurance"
    6  = "normal"
    7  = "2012-04-30"
    8  = "姚文智"
    9  = 1745
    10 = "tmp60da1"
    11 = 77
}
foreach ($col in $insRow3.Keys) {
    $wsIns.Cells.Item(3, $col).Value = $insRow3[$col]
}

# Apply header (bold/bordered) style to the newly-added header cells, and
# plain data style to the newly-added data cells, matching the existing
# columns already on the sheet.
$wsIns.Cells.Item(1, 2).Copy($wsIns.Range($wsIns.Cells.Item(1, 5), $wsIns.Cells.Item(1, 11)))
$wsIns.Cells.Item(2, 2).Copy($wsIns.Range($wsIns.Cells.Item(2, 5), $wsIns.Cells.Item(2, 11)))
$wsIns.Cells.Item(3, 2).Copy($wsIns.Range($wsIns.Cells.Item(3, 5), $wsIns.Cells.Item(3, 11)))

# Re-assert the values (copying formats with PasteSpecial-less Copy also
# copies the source value, so put the real values back afterwards).
foreach ($col in $insHeader.Keys) { $wsIns.Cells.Item(1, $col).Value = $insHeader[$col] }
foreach ($col in $insRow2.Keys)  { $wsIns.Cells.Item(2, $col).Value = $insRow2[$col] }
foreach ($col in $insRow3.Keys)  { $wsIns.Cells.Item(3, $col).Value = $insRow3[$col] }

# ---------------------------------------------------------------------------
# Sheet "債務" (debt) -> Worksheets.Item(7)
# ---------------------------------------------------------------------------
$wsDebt = $wb.Worksheets.Item(7)

$debtHeader = @{
    2  = "species"
    3  = "debtor"
    4  = "owner"
    5  = "total"
    6  = "register_date"
    7  = "register_reason"
    8  = "property_category"
    9  = "category"
    10 = "date"
    11 = "legislator_name"
    12 = "legislator_id"
    13 = "source_file"
    14 = "index"
}
foreach ($col in $debtHeader.Keys) { $wsDebt.Cells.Item(1, $col).Value = $debtHeader[$col] }

$debtRow2 = @{
    2  = "房屋貸款"
    3  = "潘瓊琪"
    4  = "花旗(台灣)商業銀行"
    5  = 8700000
    6  = "101年"
    7  = "轉貸自合庫"
    8  = "debt"
    9  = "normal"
    10 = "2012-04-30"
    11 = "姚文智"
    12 = 1745
    13 = "tmp60da1"
    14 = 87
}
foreach ($col in $debtRow2.Keys) { $wsDebt.Cells.Item(2, $col).Value = $debtRow2[$col] }

$wsDebt.Cells.Item(1, 2).Copy($wsDebt.Range($wsDebt.Cells.Item(1, 8), $wsDebt.Cells.Item(1, 14)))
$wsDebt.Cells.Item(2, 2).Copy($wsDebt.Range($wsDebt.Cells.Item(2, 8), $wsDebt.Cells.Item(2, 14)))

foreach ($col in $debtHeader.Keys) { $wsDebt.Cells.Item(1, $col).Value = $debtHeader[$col] }
foreach ($col in $debtRow2.Keys)  { $wsDebt.Cells.Item(2, $col).Value = $debtRow2[$col] }

# ---------------------------------------------------------------------------
# Sheet "事業投資" (business investment) -> Worksheets.Item(8)
# ---------------------------------------------------------------------------
$wsInv = $wb.Worksheets.Item(8)

$invHeader = @{
    2  = "owner"
    3  = "company"
    4  = "address"
    5  = "total"
    6  = "register_date"
    7  = "register_reason"
    8  = "property_category"
    9  = "category"
    10 = "date"
    11 = "legislator_name"
    12 = "legislator_id"
    13 = "source_file"
    14 = "index"
}
foreach ($col in $invHeader.Keys) { $wsInv.Cells.Item(1, $col).Value = $invHeader[$col] }

$invRow2 = @{
    2  = "潘瓊琪"
    3  = "創意企業有限公司"
    4  = "新北市新店區民權路88號4F"
    5  = 2500000
    6  = "97年"
    7  = "投資"
    8  = "investment"
    9  = "normal"
    10 = "2012-04-30"
    11 = "姚文智"
    12 = 1745
    13 = "tmp60da1"
    14 = 92
}
foreach ($col in $invRow2.Keys) { $wsInv.Cells.Item(2, $col).Value = $invRow2[$col] }

$wsInv.Cells.Item(1, 2).Copy($wsInv.Range($wsInv.Cells.Item(1, 8), $wsInv.Cells.Item(1, 14)))
$wsInv.Cells.Item(2, 2).Copy($wsInv.Range($wsInv.Cells.Item(2, 8), $wsInv.Cells.Item(2, 14)))

foreach ($col in $invHeader.Keys) { $wsInv.Cells.Item(1, $col).Value = $invHeader[$col] }
foreach ($col in $invRow2.Keys)  { $wsInv.Cells.Item(2, $col).Value = $invRow2[$col] }
